# Mise à jour site
# Updates the "disponible" (F) column quantities on Feuil1, re-formats that
# column as whole numbers, switches the vertical alignment of the detail
# rows to "top" and leaves the sheet scrolled/selected near the bottom,
# matching the author's last interactive session.

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Feuil1")

# ---------------------------------------------------------------------
# 1. New values for column F ("disponible"), rows 2-72
# ---------------------------------------------------------------------
$values = @{
    2  = 19
    3  = 0
    4  = 0
    5  = 5
    6  = 15
    7  = 8
    8  = 10
    9  = 12
    10 = 0
    11 = 26
    12 = 9
    13 = 10
    14 = 0
    15 = 12
    16 = 8
    17 = 13
    18 = 5
    19 = 10
    20 = 5
    21 = 5
    22 = 0
    23 = 6
    24 = 2
    25 = 1
    26 = 0
    27 = 17
    28 = 0
    29 = 15
    30 = 0
    31 = 0
    32 = 0
    33 = 14
    34 = 9
    35 = 0
    36 = 9
    37 = 14
    38 = 11
    39 = 0
    40 = 0
    41 = 25
    42 = 0
    43 = 0
    44 = 0
    45 = 0
    46 = 0
    47 = 12
    48 = 5
    49 = 0
    50 = 8
    51 = 1
    52 = 2
    53 = 0
    54 = 1
    55 = 2
    56 = 4
    57 = 0
    58 = 0
    59 = 0
    60 = 2
    61 = 6
    62 = 10
    63 = 0
    64 = 0
    65 = 0
    66 = -2
    67 = 0
    68 = 0
    69 = 9
    70 = 0
    71 = 24
    72 = 26
}

foreach ($row in $values.Keys) {
    $ws.Cells.Item($row, 6).Value = $values[$row]
}

# ---------------------------------------------------------------------
# 2. Number formatting: the whole "disponible" column (header + data)
#    is displayed as a plain integer ("0") instead of General/text.
# ---------------------------------------------------------------------
$ws.Cells.Item(1, 6).NumberFormat = "0"

for ($row = 2; $row -le 70; $row++) {
    $cell = $ws.Cells.Item($row, 6)
    $cell.NumberFormat = "0"
    $cell.VerticalAlignment = -4160   # xlTop
}

for ($row = 71; $row -le 72; $row++) {
    $ws.Cells.Item($row, 6).NumberFormat = "0"
}

# ---------------------------------------------------------------------
# 3. Restore the on-screen selection/scroll state used when the file was
#    last saved (bottom of the price list).
# ---------------------------------------------------------------------
$ws.Activate()
$excel.ActiveWindow.ScrollRow = 53
$excel.ActiveWindow.ScrollColumn = 1
$ws.Range("F73").Select()
